$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 3: "想去人数" (want-to-go count) bumped 159 -> 160
    $ws.Range("F3").Value = 160

    # Row 4: "想去人数" (want-to-go count) bumped 143 -> 145
    $ws.Range("F4").Value = 145

    # New row 5: a newly announced event, appended below the existing data.
    # Copy the formatting (bold / bordered / centered style) from A4 so the
    # new index cell A5 matches the rest of column A.
    $ws.Range("A4").Copy()
    $ws.Range("A5").PasteSpecial(-4122)

    $ws.Range("A5").Value = 4

    # Text cells that look like dates ("2024-06-01") need a leading
    # apostrophe so Excel stores them as literal text instead of silently
    # converting to a date serial number; ClearFormats afterwards drops the
    # implicit "Text" number-format style Excel applies, keeping the cell's
    # style identical to its untouched neighbours (no explicit s attribute).
    $ws.Range("B5").Value = "'2024-06-01"
    $ws.Range("B5").ClearFormats()

    $ws.Range("C5").Value = "丽水·动漫游戏展"
    $ws.Range("D5").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Range("E5").Value = "2024.06.01 10:00-06.01 17:00"
    $ws.Range("F5").Value = 1
    $ws.Range("G5").Value = 45
    $ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=84450"
    $ws.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202404/tdhb9QSW1713333412467.jpeg"
}
